# daily auto push: 2025-10-06 02:01 UTC
# Append the new day's row (row 68) to the bottom of the data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 68

# Column A holds a literal "yyyy/mm/dd" text label (not a real date), so force
# the cell to Text format before assigning the value; otherwise Excel would
# auto-convert the "2025/10/06" string into a date serial number. Reset the
# formatting back to the sheet's default afterwards so the new row matches
# the styling of the existing data rows.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2025/10/06"
$ws.Range("A" + $newRow).ClearFormats()

$ws.Range("B" + $newRow).Value = "月"
$ws.Range("C" + $newRow).Value = 8
$ws.Range("D" + $newRow).Value = 66
